$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 2083.8333
$ws.Range("I31").Value = 2083.8333
$ws.Range("K31").Value = 6251.499899999999
$ws.Range("M31").Value = -6021.499899999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 15790.242
$ws.Range("J87").Value = 15790.242
$ws.Range("L87").Value = 15790.242
$ws.Range("N87").Value = -18286.242

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 15790.242
$ws.Range("J90").Value = 15790.242
$ws.Range("L90").Value = 47370.726
$ws.Range("N90").Value = -59850.726

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 19232514
$ws.Range("I96").Value = 31251886
$ws.Range("J96").Value = 1519.8
$ws.Range("K96").Value = 93755658
$ws.Range("L96").Value = 4559.4
$ws.Range("M96").Value = -93754285
$ws.Range("N96").Value = -7305.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 24050798
$ws.Range("J129").Value = 28860862
$ws.Range("L129").Value = 86582586
$ws.Range("N129").Value = -86592586

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 81237.234
$ws.Range("I132").Value = 95725.82000000001
$ws.Range("J132").Value = 1550
$ws.Range("K132").Value = 287177.46
$ws.Range("L132").Value = 4650
$ws.Range("M132").Value = -284647.46
$ws.Range("N132").Value = -9710

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 866.1429000000001
$ws.Range("I2").Value = 707.6177
$ws.Range("J2").Value = 1539.875
$ws.Range("K2").Value = 707.6177
$ws.Range("L2").Value = 1539.875
$ws.Range("M2").Value = -594.6177
$ws.Range("N2").Value = -1765.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2341.8262
$ws.Range("I61").Value = 2097.4666
$ws.Range("J61").Value = 2800
$ws.Range("K61").Value = 2097.4666
$ws.Range("L61").Value = 2800
$ws.Range("M61").Value = -1885.4666
$ws.Range("N61").Value = -3224

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6903.3184
$ws.Range("I74").Value = 10421.583
$ws.Range("J74").Value = 2681.4
$ws.Range("K74").Value = 10421.583
$ws.Range("L74").Value = 2681.4
$ws.Range("M74").Value = -9547.583000000001
$ws.Range("N74").Value = -4429.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 6903.3184
$ws.Range("I77").Value = 10421.583
$ws.Range("J77").Value = 2681.4
$ws.Range("K77").Value = 52107.915
$ws.Range("L77").Value = 13407
$ws.Range("M77").Value = -47739.915
$ws.Range("N77").Value = -22143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 866.1429000000001
$ws.Range("I116").Value = 707.6177
$ws.Range("J116").Value = 1539.875
$ws.Range("K116").Value = 707.6177
$ws.Range("L116").Value = 1539.875
$ws.Range("M116").Value = 1586.3823
$ws.Range("N116").Value = -6127.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2311.3076
$ws.Range("I132").Value = 1377.7142
$ws.Range("J132").Value = 3400.5
$ws.Range("K132").Value = 4133.142599999999
$ws.Range("L132").Value = 10201.5
$ws.Range("M132").Value = -1603.142599999999
$ws.Range("N132").Value = -15261.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 85000
$ws.Range("J135").Value = 85000
$ws.Range("L135").Value = 85000
$ws.Range("N135").Value = -95140

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2341.8262
$ws.Range("I136").Value = 2097.4666
$ws.Range("J136").Value = 2800
$ws.Range("K136").Value = 6292.399800000001
$ws.Range("L136").Value = 8400
$ws.Range("M136").Value = -3742.399800000001
$ws.Range("N136").Value = -13500

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 866.1429000000001
$ws.Range("I3").Value = 707.6177
$ws.Range("J3").Value = 1539.875
$ws.Range("K3").Value = 707.6177
$ws.Range("L3").Value = 1539.875
$ws.Range("M3").Value = -593.6177
$ws.Range("N3").Value = -1767.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1051.8334
$ws.Range("I107").Value = 1002.2
$ws.Range("K107").Value = 1002.2
$ws.Range("M107").Value = 917.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 43500
$ws.Range("J138").Value = 43500
$ws.Range("L138").Value = 43500
$ws.Range("N138").Value = -53780

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1784.87
$ws.Range("I31").Value = 926.4559
$ws.Range("J31").Value = 3609
$ws.Range("K31").Value = 926.4559
$ws.Range("L31").Value = 3609
$ws.Range("M31").Value = -631.4559
$ws.Range("N31").Value = -4199

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1784.87
$ws.Range("I34").Value = 926.4559
$ws.Range("J34").Value = 3609
$ws.Range("K34").Value = 926.4559
$ws.Range("L34").Value = 3609
$ws.Range("M34").Value = -724.4559
$ws.Range("N34").Value = -4013

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2568.6562
$ws.Range("I132").Value = 1991.4584
$ws.Range("J132").Value = 4300.25
$ws.Range("K132").Value = 5974.3752
$ws.Range("L132").Value = 12900.75
$ws.Range("M132").Value = -3444.3752
$ws.Range("N132").Value = -17960.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5000
$ws.Range("I3").Value = 5000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 15000
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -14888
$ws.Range("N3").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3000
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 3090.9092
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 9272.7276
$ws.Range("M70").Value = -5685
$ws.Range("N70").Value = -9902.7276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 3000
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 3090.9092
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 9272.7276
$ws.Range("M73").Value = -4908
$ws.Range("N73").Value = -11456.7276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 10243.75
$ws.Range("I87").Value = 6991.6665
$ws.Range("K87").Value = 20974.9995
$ws.Range("M87").Value = -19726.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 10243.75
$ws.Range("I90").Value = 6991.6665
$ws.Range("K90").Value = 62924.9985
$ws.Range("M90").Value = -56684.9985

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1145.25
$ws.Range("I113").Value = 959.2917
$ws.Range("J113").Value = 1703.125
$ws.Range("K113").Value = 2877.8751
$ws.Range("L113").Value = 5109.375
$ws.Range("M113").Value = -707.8751000000002
$ws.Range("N113").Value = -9449.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2519.1035
$ws.Range("I132").Value = 2007.8
$ws.Range("J132").Value = 3655.3333
$ws.Range("K132").Value = 6023.4
$ws.Range("L132").Value = 10965.9999
$ws.Range("M132").Value = -3493.4
$ws.Range("N132").Value = -16025.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 64019.062
$ws.Range("J122").Value = 2861
$ws.Range("L122").Value = 8583
$ws.Range("N122").Value = -13483

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 267035.4
$ws.Range("I136").Value = 589001.25
$ws.Range("J136").Value = 6396.4287
$ws.Range("K136").Value = 1767003.75
$ws.Range("L136").Value = 19189.2861
$ws.Range("M136").Value = -1764453.75
$ws.Range("N136").Value = -24289.2861
